$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (C column) date for rows 2-11 from 2023-09-06 (45175) to 2023-09-14 (45183)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
